# CryCompanywiseStockReport_1.xlsx — "Add file from private repo"
#
# The underlying data rows of the stock report got shuffled: for 24 pairs of
# adjacent rows, the per-item figures (Item code, Item name, Rate, Stock
# value/rate, Qty, Value) ended up swapped between the two rows while the
# serial-number column (A) stayed anchored to its row. This script swaps
# columns B..G back between each such row pair to reproduce that edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowValues {
    param($ws, $r1, $r2, $cols)
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

$cols = @("B","C","D","E","F","G")

# Row pairs affected by the edit (1-based worksheet row numbers).
$pairs = @(
    @(142,143),
    @(256,257),
    @(305,306),
    @(308,310),
    @(338,339),
    @(342,344),
    @(347,348),
    @(364,365),
    @(367,368),
    @(371,372),
    @(381,382),
    @(392,393),
    @(411,412),
    @(413,414),
    @(423,424),
    @(528,529),
    @(571,572),
    @(573,574),
    @(575,576),
    @(578,579),
    @(585,586),
    @(679,680),
    @(701,702),
    @(707,708)
)

foreach ($pair in $pairs) {
    Swap-RowValues $ws $pair[0] $pair[1] $cols
}

Write-Output "Swapped B:G between $($pairs.Count) row pairs."
